{"js": "// Replace the two-digit-division expressions in the quiz table with the\n// newly generated set (author's commit: \"Update master to output generated\n// at c986bee\"). Each old expression occurs exactly once in the document,\n// so a literal, case-sensitive search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"55\u00f74=\", \"43\u00f77=\"],\n  [\"97\u00f75=\", \"70\u00f75=\"],\n  [\"82\u00f74=\", \"64\u00f77=\"],\n  [\"89\u00f73=\", \"10\u00f73=\"],\n  [\"62\u00f74=\", \"59\u00f79=\"],\n  [\"84\u00f76=\", \"30\u00f74=\"],\n  [\"22\u00f79=\", \"34\u00f75=\"],\n  [\"31\u00f79=\", \"43\u00f77=\"],\n  [\"11\u00f77=\", \"67\u00f75=\"],\n  [\"82\u00f72=\", \"75\u00f72=\"],\n  [\"54\u00f78=\", \"47\u00f74=\"],\n  [\"22\u00f73=\", \"15\u00f72=\"],\n  [\"86\u00f78=\", \"81\u00f75=\"],\n  [\"46\u00f72=\", \"13\u00f73=\"],\n  [\"57\u00f73=\", \"26\u00f76=\"],\n  [\"40\u00f77=\", \"16\u00f72=\"],\n  [\"71\u00f74=\", \"29\u00f74=\"],\n  [\"52\u00f75=\", \"81\u00f74=\"],\n  [\"13\u00f79=\", \"95\u00f77=\"],\n  [\"23\u00f77=\", \"71\u00f75=\"],\n  [\"98\u00f78=\", \"41\u00f77=\"],\n  [\"34\u00f79=\", \"48\u00f72=\"],\n  [\"71\u00f72=\", \"37\u00f76=\"],\n  [\"43\u00f73=\", \"29\u00f78=\"],\n  [\"74\u00f78=\", \"62\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-division expressions in the quiz table with the\n# newly generated set (author's commit: \"Update master to output generated\n# at c986bee\"). Each old expression occurs exactly once in the document,\n# so a literal, case-sensitive Find/Replace per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"55\u00f74=\", \"43\u00f77=\"),\n    @(\"97\u00f75=\", \"70\u00f75=\"),\n    @(\"82\u00f74=\", \"64\u00f77=\"),\n    @(\"89\u00f73=\", \"10\u00f73=\"),\n    @(\"62\u00f74=\", \"59\u00f79=\"),\n    @(\"84\u00f76=\", \"30\u00f74=\"),\n    @(\"22\u00f79=\", \"34\u00f75=\"),\n    @(\"31\u00f79=\", \"43\u00f77=\"),\n    @(\"11\u00f77=\", \"67\u00f75=\"),\n    @(\"82\u00f72=\", \"75\u00f72=\"),\n    @(\"54\u00f78=\", \"47\u00f74=\"),\n    @(\"22\u00f73=\", \"15\u00f72=\"),\n    @(\"86\u00f78=\", \"81\u00f75=\"),\n    @(\"46\u00f72=\", \"13\u00f73=\"),\n    @(\"57\u00f73=\", \"26\u00f76=\"),\n    @(\"40\u00f77=\", \"16\u00f72=\"),\n    @(\"71\u00f74=\", \"29\u00f74=\"),\n    @(\"52\u00f75=\", \"81\u00f74=\"),\n    @(\"13\u00f79=\", \"95\u00f77=\"),\n    @(\"23\u00f77=\", \"71\u00f75=\"),\n    @(\"98\u00f78=\", \"41\u00f77=\"),\n    @(\"34\u00f79=\", \"48\u00f72=\"),\n    @(\"71\u00f72=\", \"37\u00f76=\"),\n    @(\"43\u00f73=\", \"29\u00f78=\"),\n    @(\"74\u00f78=\", \"62\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
